$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "svin_gylle" data row (row 10) entirely - the rows below shift up.
$ws.Rows.Item(10).Delete()

# Remove the "kvæg_gylle" data row, which after the previous delete is now row 15.
$ws.Rows.Item(15).Delete()

# Zero out a handful of cells that changed independently of the row shuffle.
# Row 3 = spalter_smågrise
$ws.Range("E3").Value = 0
$ws.Range("I3:L3").Value = 0

# Row 9 = farestald_fuldspalte
$ws.Range("E9").Value = 0
$ws.Range("I9:L9").Value = 0

# Leave the cursor where the edit session left it.
$ws.Range("L8").Select()
